$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.8196770910878
$ws.Range("C2").Value = 8.701991038463163
$ws.Range("D2").Value = 7.740546356698469
$ws.Range("E2").Value = 13.1147000458921
$ws.Range("F2").Value = 39.80291378266475
$ws.Range("J2").Value = 10.43965123685456
$ws.Range("K2").Value = 10.97304233228034
$ws.Range("L2").Value = 10.64414207485471
$ws.Range("O2").Value = 30.85463545647013

$ws.Range("B3").Value = 14.62931646709673
$ws.Range("C3").Value = 8.691439001280019
$ws.Range("D3").Value = 7.723530605890372
$ws.Range("E3").Value = 13.13361988409823
$ws.Range("F3").Value = 39.89908891827489
$ws.Range("J3").Value = 10.46108441887647
$ws.Range("K3").Value = 10.83240246286649
$ws.Range("L3").Value = 10.64266630505112
$ws.Range("O3").Value = 30.94776686950024

$ws.Range("B4").Value = 14.513768673553
$ws.Range("C4").Value = 8.685101927765675
$ws.Range("D4").Value = 7.714087536435531
$ws.Range("E4").Value = 13.14684489049918
$ws.Range("F4").Value = 39.96572229615035
$ws.Range("J4").Value = 10.47505103561553
$ws.Range("K4").Value = 10.746781596787
$ws.Range("L4").Value = 10.64300740739907
$ws.Range("O4").Value = 31.01030512843756

$ws.Range("B5").Value = 14.46707100463819
$ws.Range("C5").Value = 8.682555940822406
$ws.Range("D5").Value = 7.710494774059507
$ws.Range("E5").Value = 13.15263905169115
$ws.Range("F5").Value = 39.99478004119068
$ws.Range("J5").Value = 10.48094584146439
$ws.Range("K5").Value = 10.71211283334718
$ws.Range("L5").Value = 10.64346122514592
$ws.Range("O5").Value = 31.03713545824397

$ws.Range("B6").Value = 14.45934193830342
$ws.Range("C6").Value = 8.682135402780364
$ws.Range("D6").Value = 7.709913698467336
$ws.Range("E6").Value = 13.15362563393744
$ws.Range("F6").Value = 39.99971998852507
$ws.Range("J6").Value = 10.48193696224501
$ws.Range("K6").Value = 10.7063706581296
$ws.Range("L6").Value = 10.64355562838779
$ws.Range("O6").Value = 31.04167184409202

$ws.Range("B7").Value = 14.51313724720199
$ws.Range("C7").Value = 8.685067443257244
$ws.Range("D7").Value = 7.714038045733394
$ws.Range("E7").Value = 13.14692139268847
$ws.Range("F7").Value = 39.96610647267243
$ws.Range("J7").Value = 10.47512971125002
$ws.Range("K7").Value = 10.74631309077753
$ws.Range("L7").Value = 10.64301225155665
$ws.Range("O7").Value = 31.01066152538208

$ws.Range("B8").Value = 14.75379599158112
$ws.Range("C8").Value = 8.698323681061934
$ws.Range("D8").Value = 7.73447249540577
$ws.Range("E8").Value = 13.12089013287786
$ws.Range("F8").Value = 39.83450011184441
$ws.Range("J8").Value = 10.44687428678329
$ws.Range("K8").Value = 10.9244206438578
$ws.Range("L8").Value = 10.64337510359837
$ws.Range("O8").Value = 30.88563503520292

$ws.Range("B9").Value = 15.23395264602096
$ws.Range("C9").Value = 8.725413719841216
$ws.Range("D9").Value = 7.782389371023143
$ws.Range("E9").Value = 13.08258103744324
$ws.Range("F9").Value = 39.63667211046047
$ws.Range("J9").Value = 10.39784338144481
$ws.Range("K9").Value = 11.27781867224489
$ws.Range("L9").Value = 10.65392819927284
$ws.Range("O9").Value = 30.68299857646537

$ws.Range("B10").Value = 15.58858143720893
$ws.Range("C10").Value = 8.74594301090176
$ws.Range("D10").Value = 7.822197323066228
$ws.Range("E10").Value = 13.06217048981767
$ws.Range("F10").Value = 39.52817929400722
$ws.Range("J10").Value = 10.36567830589268
$ws.Range("K10").Value = 11.53773602245055
$ws.Range("L10").Value = 10.6676042628551
$ws.Range("O10").Value = 30.56012022536088

$ws.Range("B11").Value = 15.74965573412391
$ws.Range("C11").Value = 8.75540948033442
$ws.Range("D11").Value = 7.841265088908213
$ws.Range("E11").Value = 13.05455763580059
$ws.Range("F11").Value = 39.48684563244928
$ws.Range("J11").Value = 10.35187691088274
$ws.Range("K11").Value = 11.65557283967936
$ws.Range("L11").Value = 10.67509371608225
$ws.Range("O11").Value = 30.50987806890969

$ws.Range("B12").Value = 15.81056100720925
$ws.Range("C12").Value = 8.759011739754985
$ws.Range("D12").Value = 7.848619669007608
$ws.Range("E12").Value = 13.05191459984163
$ws.Range("F12").Value = 39.4723481451374
$ws.Range("J12").Value = 10.34676964745901
$ws.Range("K12").Value = 11.70009901126914
$ws.Range("L12").Value = 10.67811041067323
$ws.Range("O12").Value = 30.4916669984516

$ws.Range("B13").Value = 15.79744882310344
$ws.Range("C13").Value = 8.758235163220933
$ws.Range("D13").Value = 7.847029830575998
$ws.Range("E13").Value = 13.0524731711102
$ws.Range("F13").Value = 39.47541906554273
$ws.Range("J13").Value = 10.34786430101282
$ws.Range("K13").Value = 11.69051437676893
$ws.Range("L13").Value = 10.67745270820898
$ws.Range("O13").Value = 30.49555283219104

$ws.Range("B14").Value = 15.75466855903182
$ws.Range("C14").Value = 8.755705494209323
$ws.Range("D14").Value = 7.841867491098363
$ws.Range("E14").Value = 13.05433538981734
$ws.Range("F14").Value = 39.48562976733918
$ws.Range("J14").Value = 10.35145435020541
$ws.Range("K14").Value = 11.65923817634425
$ws.Range("L14").Value = 10.67533829224253
$ws.Range("O14").Value = 30.50836350163537

$ws.Range("B15").Value = 15.72845109311726
$ws.Range("C15").Value = 8.754158254684805
$ws.Range("D15").Value = 7.838722746782634
$ws.Range("E15").Value = 13.05550725991862
$ws.Range("F15").Value = 39.49203452487178
$ws.Range("J15").Value = 10.35366884715005
$ws.Range("K15").Value = 11.64006696670724
$ws.Range("L15").Value = 10.67406661726975
$ws.Range("O15").Value = 30.51631651994669

$ws.Range("B16").Value = 15.5780451788362
$ws.Range("C16").Value = 8.745326851171599
$ws.Range("D16").Value = 7.820970169095691
$ws.Range("E16").Value = 13.06270159857967
$ws.Range("F16").Value = 39.53104205833136
$ws.Range("J16").Value = 10.36659693828429
$ws.Range("K16").Value = 11.5300237336958
$ws.Range("L16").Value = 10.66714017980735
$ws.Range("O16").Value = 30.56351758589214

$ws.Range("B17").Value = 15.48567333709893
$ws.Range("C17").Value = 8.739941159232622
$ws.Range("D17").Value = 7.810322442130212
$ws.Range("E17").Value = 13.06754288780133
$ws.Range("F17").Value = 39.55702717098796
$ws.Range("J17").Value = 10.37474035316427
$ws.Range("K17").Value = 11.46238543053433
$ws.Range("L17").Value = 10.66321461743771
$ws.Range("O17").Value = 30.59392340894292

$ws.Range("B18").Value = 15.43252331576831
$ws.Range("C18").Value = 8.736855511680414
$ws.Range("D18").Value = 7.804288650899793
$ws.Range("E18").Value = 13.07048488790463
$ws.Range("F18").Value = 39.57272796320003
$ws.Range("J18").Value = 10.37950244055532
$ws.Range("K18").Value = 11.42344617205864
$ws.Range("L18").Value = 10.6610762125933
$ws.Range("O18").Value = 30.61194431663885

$ws.Range("B19").Value = 15.41452584573247
$ws.Range("C19").Value = 8.73581286105358
$ws.Range("D19").Value = 7.802261371495247
$ws.Range("E19").Value = 13.07150805433448
$ws.Range("F19").Value = 39.57817357987967
$ws.Range("J19").Value = 10.38112824720555
$ws.Range("K19").Value = 11.41025709967097
$ws.Range("L19").Value = 10.660372757315
$ws.Range("O19").Value = 30.61813726739934

$ws.Range("B20").Value = 15.49550896449624
$ws.Range("C20").Value = 8.740513233190844
$ws.Range("D20").Value = 7.811446571974791
$ws.Range("E20").Value = 13.06701123660785
$ws.Range("F20").Value = 39.55418287707878
$ws.Range("J20").Value = 10.37386538132552
$ws.Range("K20").Value = 11.46958959871588
$ws.Range("L20").Value = 10.66362014748232
$ws.Range("O20").Value = 30.59063155868929

$ws.Range("B21").Value = 15.7672370302254
$ws.Range("C21").Value = 8.75644805071498
$ws.Range("D21").Value = 7.843380189280754
$ws.Range("E21").Value = 13.05378190867001
$ws.Range("F21").Value = 39.48259928984561
$ws.Range("J21").Value = 10.35039663876399
$ws.Range("K21").Value = 11.66842766286213
$ws.Range("L21").Value = 10.67595446047003
$ws.Range("O21").Value = 30.50457858143051

$ws.Range("B22").Value = 15.94427789621938
$ws.Range("C22").Value = 8.766964060888302
$ws.Range("D22").Value = 7.865029996599124
$ws.Range("E22").Value = 13.04653310692512
$ws.Range("F22").Value = 39.44254600982502
$ws.Range("J22").Value = 10.33575203679623
$ws.Range("K22").Value = 11.79780280967133
$ws.Range("L22").Value = 10.68506744734968
$ws.Range("O22").Value = 30.45308629831329

$ws.Range("B23").Value = 15.84985568064147
$ws.Range("C23").Value = 8.761342441895803
$ws.Range("D23").Value = 7.853405090310261
$ws.Range("E23").Value = 13.05027429805415
$ws.Range("F23").Value = 39.4633069655069
$ws.Range("J23").Value = 10.3435048124319
$ws.Range("K23").Value = 11.72881806572068
$ws.Range("L23").Value = 10.68010803968909
$ws.Range("O23").Value = 30.48013383520353

$ws.Range("B24").Value = 15.49106241092084
$ws.Range("C24").Value = 8.740254565294087
$ws.Range("D24").Value = 7.810938079328689
$ws.Range("E24").Value = 13.06725110180393
$ws.Range("F24").Value = 39.55546641009427
$ws.Range("J24").Value = 10.37426070584053
$ws.Range("K24").Value = 11.46633275672636
$ws.Range("L24").Value = 10.66343643825322
$ws.Range("O24").Value = 30.59211812158565

$ws.Range("B25").Value = 15.10350454338566
$ws.Range("C25").Value = 8.717972166717479
$ws.Range("D25").Value = 7.7686045438898
$ws.Range("E25").Value = 13.09158393568018
$ws.Range("F25").Value = 39.68372597824289
$ws.Range("J25").Value = 10.4104278822533
$ws.Range("K25").Value = 11.18200601018628
$ws.Range("L25").Value = 10.65002727404123
$ws.Range("O25").Value = 30.73325596472768
